$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.169.07'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '4.037.65'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '539.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '152.33'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').Value = '4.033.57'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '53.54'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +10.45%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000329'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '4.680.64'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '4.052.79'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.61'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('D21').Value = '72.122.69'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '447.99'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.61%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '97.65'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.50'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '14.62'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '4.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +17.95%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.28'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.95'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '37.15'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.15'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +18.02%  '
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '13.59'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '49.22'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +14.69%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '679.71'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '66.91'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.456'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.19%  '
$ws.Range('D39').Value = '0.0₃0880'
$ws.Range('E39').Value = '  +3.54%  '
$ws.Range('E40').Value = '  -5.79%  '
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.22'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +16.38%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.37'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.75%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.63'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('E51').Value = '  +1.28%  '
